# "switching to # because of excel rules"
# The Main sheet keeps a small lookup table in column A/B(/C) whose values
# are the names of other worksheets in the workbook (DefaultCoefficients,
# DebugCoefficients, CalibrationMatrix, StructArray, AnArray, Thresholds).
# Those reference values are being re-pointed to a "#"-prefixed alias of
# each sheet name (Excel sheet-naming rules forbid literal "#" in a sheet
# name, so the lookup values get the prefix instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("B25").Value = "#Thresholds"
$ws.Range("B18").Value = "#AnArray"
$ws.Range("B16").Value = "#StructArray"
$ws.Range("B14").Value = "#CalibrationMatrix"
$ws.Range("B13").Value = "#DefaultCoefficients"
$ws.Range("C13").Value = "#DebugCoefficients"

# Update the active selection on the Main sheet to match the saved view.
$ws.Activate()
$ws.Range("F14").Select()
